# Refresh the "cryptos" price table (GitHub Actions data pull).
#
# Every data cell in D2:E51 (and a couple of B/C cells for two rows whose
# ranking order flipped) is plain text in the workbook ("t=inlineStr" /
# shared-string text) -- never a real number. Several of the new Price
# values (column D) are strings that Excel's COM layer would otherwise
# auto-coerce into numbers on a plain `.Value = "..."` assignment (e.g.
# "0.999", "11.10", "0.0350" would become 0.999, 11.1, 0.035 and lose
# their text type / formatting). To keep those cells genuinely text
# (matching the original file) without forcing a NumberFormat change
# (which would add a style index that isn't present in the target), we
# write them as a `="..."` text formula first and then flatten the whole
# Price column back down to literal values via Copy/PasteSpecial(values
# only) -- the same thing a user gets from Paste Values in the UI.
#
# Values that can never be mis-parsed as numbers (two decimal points,
# coin names, URLs, the padded "  +x.xx%  " volume strings) are just
# assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.493.51"
$ws.Range("E2").Value = "  +4.38%  "
$ws.Range("D3").Value = "2.972.67"
$ws.Range("E3").Value = "  +2.70%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Formula = "=""579.87"""
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("D6").Formula = "=""152.77"""
$ws.Range("E6").Value = "  +6.51%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "2.970.49"
$ws.Range("E8").Value = "  +2.90%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").Formula = "=""0.511"""
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("D10").Formula = "=""6.95"""
$ws.Range("E10").Value = "  +4.99%  "
$ws.Range("D11").Formula = "=""0.153"""
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("D12").Formula = "=""0.446"""
$ws.Range("E12").Value = "  +3.18%  "
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("D14").Formula = "=""34.37"""
$ws.Range("E14").Value = "  +7.63%  "
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").Value = "3.467.83"
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("D17").Value = "64.466.51"
$ws.Range("E17").Value = "  +4.46%  "
$ws.Range("D18").Formula = "=""6.92"""
$ws.Range("E18").Value = "  +4.80%  "
$ws.Range("D19").Value = "2.964.18"
$ws.Range("E19").Value = "  +2.47%  "
$ws.Range("D20").Formula = "=""445.59"""
$ws.Range("E20").Value = "  +1.97%  "
$ws.Range("D21").Formula = "=""13.61"""
$ws.Range("E21").Value = "  +3.31%  "
$ws.Range("D22").Formula = "=""0.677"""
$ws.Range("E22").Value = "  +3.70%  "
$ws.Range("D23").Formula = "=""7.25"""
$ws.Range("E23").Value = "  +5.29%  "
$ws.Range("D24").Formula = "=""80.67"""
$ws.Range("E24").Value = "  +2.06%  "
$ws.Range("D25").Formula = "=""11.10"""
$ws.Range("E25").Value = "  +9.89%  "
$ws.Range("D26").Formula = "=""12.24"""
$ws.Range("E26").Value = "  +2.82%  "
$ws.Range("E27").Value = "  +8.14%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Formula = "=""7.79"""
$ws.Range("E29").Value = "  +10.51%  "
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("D31").Formula = "=""2.24"""
$ws.Range("E31").Value = "  +9.09%  "
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("D33").Formula = "=""0.110"""
$ws.Range("E33").Value = "  +3.35%  "
$ws.Range("D34").Formula = "=""26.50"""
$ws.Range("E34").Value = "  +3.57%  "
$ws.Range("D35").Formula = "=""0.999"""
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Formula = "=""0.975"""
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("D37").Formula = "=""5.65"""
$ws.Range("E37").Value = "  +3.98%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Formula = "=""2.10"""
$ws.Range("E38").Value = "  +8.12%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Formula = "=""3.06"""
$ws.Range("E39").Value = "  +5.72%  "
$ws.Range("D40").Formula = "=""48.55"""
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").Formula = "=""44.47"""
$ws.Range("E41").Value = "  +13.22%  "
$ws.Range("D42").Formula = "=""0.120"""
$ws.Range("E42").Value = "  +2.96%  "
$ws.Range("E43").Value = "  +9.22%  "
$ws.Range("D44").Formula = "=""8.36"""
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("D45").Formula = "=""387.59"""
$ws.Range("E45").Value = "  +14.86%  "
$ws.Range("D46").Value = "2.777.28"
$ws.Range("E46").Value = "  +3.42%  "
$ws.Range("D47").Formula = "=""0.0350"""
$ws.Range("E47").Value = "  +5.42%  "
$ws.Range("D48").Formula = "=""135.17"""
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Formula = "=""0.000223"""
$ws.Range("E50").Value = "  +14.07%  "
$ws.Range("E51").Value = "  +2.55%  "

# Flatten the text-formula cells in the Price column back into literal
# string values (Copy + Paste Values), so no formula/"f" element remains
# and no cell style/number-format changes.
$ws.Range("D2:D51").Copy()
$ws.Range("D2:D51").PasteSpecial(-4163)
$excel.CutCopyMode = 0
